# Applies updated TPM-based NATMI values to the Ifnb1-Ifnar2 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new numeric value, as produced by the refreshed TPM pipeline.
$updates = @{
    "I2" = 0.8950246002264373
    "J2" = 0.8950246002264374
    "M2" = 30.7693535
    "N2" = 61.538707
    "O2" = 0.1179731387858698
    "P2" = 0.08351747770158975
    "Q2" = 51.93988922568883
    "R2" = 311.639335354133
    "S2" = 0.1055888613792811
    "T2" = 0.07475019709178576
    "I3" = 0.8950246002264373
    "J3" = 0.8950246002264374
    "O3" = 0.09273042782012855
    "P3" = 0.09847086613229204
    "S3" = 0.08299601408853706
    "T3" = 0.08813384759400572
    "I4" = 0.8950246002264373
    "J4" = 0.8950246002264374
    "M4" = 81.77185533333333
    "N4" = 245.315566
    "O4" = 0.3135224286729781
    "P4" = 0.3329309033622996
    "Q4" = 138.0341354195949
    "R4" = 1242.307218776354
    "S4" = 0.280610286385054
    "T4" = 0.2979813486848688
    "I5" = 0.8950246002264373
    "J5" = 0.8950246002264374
    "M5" = 14.8441875
    "N5" = 29.688375
    "O5" = 0.05691427322647431
    "P5" = 0.04029168498874919
    "Q5" = 25.0575773194375
    "R5" = 150.345463916625
    "S5" = 0.0509396746417034
    "T5" = 0.03606204924950479
    "I6" = 0.8950246002264373
    "J6" = 0.8950246002264374
    "M6" = 53.27148833333334
    "N6" = 159.814465
    "O6" = 0.2042488376129897
    "P6" = 0.2168927763956593
    "Q6" = 89.92438540903723
    "R6" = 809.319468681335
    "S6" = 0.1828077342312806
    "T6" = 0.1941243704855271
    "I7" = 0.8950246002264373
    "J7" = 0.8950246002264374
    "M7" = 55.97408466666666
    "N7" = 167.922254
    "O7" = 0.2146108938815595
    "P7" = 0.22789629141941
    "Q7" = 94.48647522269177
    "R7" = 850.3782770042259
    "S7" = 0.1920820295005812
    "T7" = 0.2039727871207451
    "E8" = 2
    "F8" = 0.6666666666666666
    "G8" = 0.175453
    "H8" = 0.526359
    "I8" = 0.0930278797853264
    "J8" = 0.09302787978532641
    "M8" = 30.7693535
    "N8" = 61.538707
    "O8" = 0.1179731387858698
    "P8" = 0.08351747770158975
    "Q8" = 5.3985753796355
    "R8" = 32.391452277813
    "S8" = 0.01097479097286952
    "T8" = 0.00776945387559717
    "E9" = 2
    "F9" = 0.6666666666666666
    "G9" = 0.175453
    "H9" = 0.526359
    "I9" = 0.0930278797853264
    "J9" = 0.09302787978532641
    "O9" = 0.09273042782012855
    "P9" = 0.09847086613229204
    "Q9" = 4.243442276139334
    "R9" = 38.190980485254
    "S9" = 0.008626515091692806
    "T9" = 0.009160535896911834
    "E10" = 2
    "F10" = 0.6666666666666666
    "G10" = 0.175453
    "H10" = 0.526359
    "I10" = 0.0930278797853264
    "J10" = 0.09302787978532641
    "M10" = 81.77185533333333
    "N10" = 245.315566
    "O10" = 0.3135224286729781
    "P10" = 0.3329309033622996
    "Q10" = 14.34711733379933
    "R10" = 129.124056004194
    "S10" = 0.02916632680459338
    "T10" = 0.03097185605480813
    "E11" = 2
    "F11" = 0.6666666666666666
    "G11" = 0.175453
    "H11" = 0.526359
    "I11" = 0.0930278797853264
    "J11" = 0.09302787978532641
    "M11" = 14.8441875
    "N11" = 29.688375
    "O11" = 0.05691427322647431
    "P11" = 0.04029168498874919
    "Q11" = 2.6044572294375
    "R11" = 15.626743376625
    "S11" = 0.005294614167781673
    "T11" = 0.0037482500274816
    "E12" = 2
    "F12" = 0.6666666666666666
    "G12" = 0.175453
    "H12" = 0.526359
    "I12" = 0.0930278797853264
    "J12" = 0.09302787978532641
    "M12" = 53.27148833333334
    "N12" = 159.814465
    "O12" = 0.2042488376129897
    "P12" = 0.2168927763956593
    "Q12" = 9.346642442548333
    "R12" = 84.11978198293501
    "S12" = 0.01900083631175386
    "T12" = 0.02017707512884108
    "E13" = 2
    "F13" = 0.6666666666666666
    "G13" = 0.175453
    "H13" = 0.526359
    "I13" = 0.0930278797853264
    "J13" = 0.09302787978532641
    "M13" = 55.97408466666666
    "N13" = 167.922254
    "O13" = 0.2146108938815595
    "P13" = 0.22789629141941
    "Q13" = 9.820821077020666
    "R13" = 88.38738969318599
    "S13" = 0.01996479643663516
    "T13" = 0.02120070880168659
    "G14" = 0.02253333333333333
    "H14" = 0.06759999999999999
    "I14" = 0.01194751998823629
    "J14" = 0.01194751998823629
    "M14" = 30.7693535
    "N14" = 61.538707
    "O14" = 0.1179731387858698
    "P14" = 0.08351747770158975
    "Q14" = 0.6933360988666667
    "R14" = 4.1600165932
    "S14" = 0.001409486433719153
    "T14" = 0.0009978267342068221
    "G15" = 0.02253333333333333
    "H15" = 0.06759999999999999
    "I15" = 0.01194751998823629
    "J15" = 0.01194751998823629
    "O15" = 0.09273042782012855
    "P15" = 0.09847086613229204
    "Q15" = 0.5449829828444445
    "R15" = 4.904846845599999
    "S15" = 0.001107898639898688
    "T15" = 0.001176482641374499
    "G16" = 0.02253333333333333
    "H16" = 0.06759999999999999
    "I16" = 0.01194751998823629
    "J16" = 0.01194751998823629
    "M16" = 81.77185533333333
    "N16" = 245.315566
    "O16" = 0.3135224286729781
    "P16" = 0.3329309033622996
    "Q16" = 1.842592473511111
    "R16" = 16.5833322616
    "S16" = 0.003745815483330792
    "T16" = 0.003977698622622638
    "G17" = 0.02253333333333333
    "H17" = 0.06759999999999999
    "I17" = 0.01194751998823629
    "J17" = 0.01194751998823629
    "M17" = 14.8441875
    "N17" = 29.688375
    "O17" = 0.05691427322647431
    "P17" = 0.04029168498874919
    "Q17" = 0.334489025
    "R17" = 2.00693415
    "S17" = 0.0006799844169892433
    "T17" = 0.0004813857117628009
    "G18" = 0.02253333333333333
    "H18" = 0.06759999999999999
    "I18" = 0.01194751998823629
    "J18" = 0.01194751998823629
    "M18" = 53.27148833333334
    "N18" = 159.814465
    "O18" = 0.2042488376129897
    "P18" = 0.2168927763956593
    "Q18" = 1.200384203777778
    "R18" = 10.803457834
    "S18" = 0.002440267069955222
    "T18" = 0.002591330781291204
    "G19" = 0.02253333333333333
    "H19" = 0.06759999999999999
    "I19" = 0.01194751998823629
    "J19" = 0.01194751998823629
    "M19" = 55.97408466666666
    "N19" = 167.922254
    "O19" = 0.2146108938815595
    "P19" = 0.22789629141941
    "Q19" = 0.5449829828444445
    "R19" = 11.3515443704
    "S19" = 0.00256406794434319
    "T19" = 0.02722795496978324
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
